$d = $word.ActiveDocument

$pairs = @(
    @("121÷3=40, 1", "209÷9=23, 2"),
    @("691÷5=138, 1", "721÷3=240, 1"),
    @("856÷4=214, 0", "167÷3=55, 2"),
    @("522÷7=74, 4", "255÷7=36, 3"),
    @("170÷7=24, 2", "470÷2=235, 0"),
    @("961÷2=480, 1", "618÷9=68, 6"),
    @("188÷4=47, 0", "591÷9=65, 6"),
    @("229÷4=57, 1", "900÷4=225, 0"),
    @("965÷6=160, 5", "110÷8=13, 6"),
    @("583÷6=97, 1", "103÷2=51, 1"),
    @("601÷8=75, 1", "939÷2=469, 1"),
    @("345÷4=86, 1", "677÷6=112, 5"),
    @("127÷4=31, 3", "534÷3=178, 0"),
    @("886÷9=98, 4", "547÷5=109, 2"),
    @("422÷8=52, 6", "630÷9=70, 0"),
    @("760÷2=380, 0", "517÷2=258, 1"),
    @("815÷3=271, 2", "791÷4=197, 3"),
    @("249÷6=41, 3", "327÷6=54, 3"),
    @("538÷5=107, 3", "191÷2=95, 1"),
    @("896÷8=112, 0", "387÷2=193, 1"),
    @("556÷5=111, 1", "938÷5=187, 3"),
    @("127÷2=63, 1", "728÷8=91, 0"),
    @("189÷4=47, 1", "659÷5=131, 4"),
    @("558÷5=111, 3", "339÷5=67, 4"),
    @("622÷5=124, 2", "107÷4=26, 3")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
